# Insert a new header row at the top of Sheet1, fill it with labels
# ("Customer", "Vehicle", "Color", "Labor Type"), make the header bold,
# set the page orientation to portrait, and move the active selection
# to C6 - matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift all existing data down one row.
$ws.Rows.Item(1).Insert()

# Populate the new header row. The order in which the new strings are
# first written controls their position in the shared-string table, so
# write B1/C1/A1/D1 in that order to reproduce the authored string order
# (Vehicle, Color, Customer, Labor Type).
$ws.Range("B1").Value = "Vehicle"
$ws.Range("C1").Value = "Color"
$ws.Range("A1").Value = "Customer"
$ws.Range("D1").Value = "Labor Type"

# Bold the new header row.
$ws.Range("A1:D1").Font.Bold = $true

# Touch page setup so a pageSetup element is emitted for the sheet.
$ws.PageSetup.Orientation = 1

# Move the active cell/selection to C6.
$ws.Range("C6").Select()
